# Automatische test-sync: 2025-08-05 16:24:50
#
# Adds a new "Testmail #3" complaint row to the Logs sheet, updates the
# matching category tally on the Dashboard sheet, extends the conditional
# formatting ranges on Logs to cover the new row, and extends the bar
# chart's category/value series references to include the new Dashboard row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Logs sheet: append row 5 with the new test mail entry
# ---------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A5").Value = "Ik ben niet tevreden over hoe dit is gegaan."
$logs.Range("B5").Value = "mailmind.test@zohomail.eu"
$logs.Range("C5").Value = "Testmail #3: Ik ben niet tevreden over hoe dit is gegaan."
$logs.Range("D5").Value = "Klacht / Probleem"
$logs.Range("E5").Value = "Bedankt, we hebben dit doorgestuurd naar klachten@bedrijf.nl."
$logs.Range("F5").Value = "2025-08-05 16:24:14"
$logs.Range("G5").Value = "Ja"
$logs.Range("H5").Value = "Ja"
$logs.Range("I5").Value = "Nee"
$logs.Range("J5").Value = "Nee"

# ---------------------------------------------------------------------
# 2. Logs sheet: extend conditional formatting ranges from row 4 to row 5
# ---------------------------------------------------------------------
$colsToExtend = @("D", "G", "H", "I", "J")
foreach ($col in $colsToExtend) {
    $oldRange = $logs.Range($col + "2:" + $col + "4")
    $newRange = $logs.Range($col + "2:" + $col + "5")
    $fcs = $oldRange.FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}

# ---------------------------------------------------------------------
# 3. Dashboard sheet: append row 4 with the new category tally
# ---------------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Range("A4").Value = "Klacht / Probleem"
$dash.Range("B4").Value = 1

# ---------------------------------------------------------------------
# 4. Dashboard chart: extend category/value series references to row 4
# ---------------------------------------------------------------------
$chartObj = $dash.ChartObjects().Item(1)
$chart = $chartObj.Chart
$series = $chart.SeriesCollection(1)
$series.Formula = "=SERIES(Dashboard!`$B`$1,Dashboard!`$A`$2:`$A`$4,Dashboard!`$B`$2:`$B`$4,1)"
